# Apply the StructureDefinition-detected.xlsx update:
#  - Metadata sheet: bump Version, update Date, replace Publisher contact info
#    with Publisher/Jurisdiction info, and drop the now-redundant duplicate
#    "Contact" row.
#  - Elements sheet: update the root Extension row's Short/Definition text to
#    match the new Name/Description ("Detected" / "When the insight
#    evaluation occurred.").

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be "Contact" / "No display for ContactDetail" -- becomes
# "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row that
# no longer applies -- remove it entirely (rows below shift up).
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root "Extension" row's Short/Definition columns (K/L) now mirror the
# updated Name/Description.
$elements.Range("K2").Value = "Detected"
$elements.Range("L2").Value = "When the insight evaluation occurred."
